# Apply text replacements per the diff.
$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-01 Monday", "2024-07-02 Tuesday"),
    @("791×8=6328", "338×2=676"),
    @("143×2=286", "269×7=1883"),
    @("181×9=1629", "811×4=3244"),
    @("133×2=266", "747×2=1494"),
    @("333×3=999", "130×9=1170"),
    @("830×6=4980", "347×2=694"),
    @("379×2=758", "332×7=2324"),
    @("517×3=1551", "233×9=2097"),
    @("276×8=2208", "957×8=7656"),
    @("982×7=6874", "405×2=810"),
    @("167×7=1169", "947×3=2841"),
    @("216×4=864", "754×6=4524"),
    @("130×2=260", "820×3=2460"),
    @("703×8=5624", "372×5=1860"),
    @("510×6=3060", "907×7=6349"),
    @("929×7=6503", "160×9=1440"),
    @("802×2=1604", "129×4=516"),
    @("105×7=735", "372×3=1116"),
    @("412×2=824", "641×6=3846"),
    @("725×6=4350", "872×8=6976"),
    @("697×8=5576", "777×6=4662"),
    @("140×6=840", "745×9=6705"),
    @("279×6=1674", "349×4=1396"),
    @("231×9=2079", "233×8=1864"),
    @("566×6=3396", "621×6=3726")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
